$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting existing rows 20-45 down to 21-46
$ws.Rows("20:20").Insert()

# Populate the new row 20 with the new weekly record
$ws.Range("A20").Value = 8
$ws.Range("B20").Value = "Terminal La Palmera de La Serena"
$ws.Range("C20").Value = "Coquimbo"
$ws.Range("D20").Value = 44484
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = 100112052
$ws.Range("G20").Value = "Albahaca"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 840
$ws.Range("K20").Value = 3500
$ws.Range("L20").Value = 4000
$ws.Range("M20").Value = 3750
$ws.Range("N20").Value = "$/paquete"
$ws.Range("O20").Value = "Región de Arica y Parinacota"
$ws.Range("P20").Value = 3750
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Hortaliza"
